$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos table refresh (prices/volume %) for the latest GitHub Actions run.
# A leading apostrophe is used for Price values that would otherwise be
# auto-parsed as numbers by Excel, so the cells stay plain text like the rest
# of the sheet (e.g. "113.24" must stay text, not become the number 113.24).
$ws.Range("D2").Value = "43.027.46"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "2.288.15"
$ws.Range("E3").Value = "  +1.92%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'113.24"
$ws.Range("E5").Value = "  -1.57%  "
$ws.Range("D6").Value = "'310.30"
$ws.Range("E6").Value = "  +6.95%  "
$ws.Range("D7").Value = "'0.633"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").Value = "'0.614"
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("E10").Value = "  -4.36%  "
$ws.Range("D11").Value = "'0.0928"
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").Value = "'55.16"
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("E13").Value = "  -4.02%  "
$ws.Range("D14").Value = "'1.06"
$ws.Range("E14").Value = "  +19.12%  "
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "'15.51"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "2.633.16"
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").Value = "2.284.25"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("D19").Value = "43.084.99"
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "'7.21"
$ws.Range("E21").Value = "  -3.30%  "
$ws.Range("D22").Value = "'75.53"
$ws.Range("E22").Value = "  +2.23%  "
$ws.Range("D23").Value = "'3.64"
$ws.Range("E23").Value = "  +6.51%  "
$ws.Range("E24").Value = "  +4.15%  "
$ws.Range("D25").Value = "'256.85"
$ws.Range("E25").Value = "  +10.40%  "
$ws.Range("D26").Value = "'8.98"
$ws.Range("E26").Value = "  -2.86%  "
$ws.Range("D27").Value = "'11.79"
$ws.Range("E27").Value = "  -3.22%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("D30").Value = "'38.34"
$ws.Range("E30").Value = "  -4.57%  "
$ws.Range("D31").Value = "'175.27"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").Value = "'22.25"
$ws.Range("E32").Value = "  +4.57%  "
$ws.Range("E33").Value = "  -2.83%  "
$ws.Range("D34").Value = "'0.0902"
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("D35").Value = "'5.72"
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").Value = "'5.03"
$ws.Range("E36").Value = "  +6.97%  "
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("D38").Value = "'4.21"
$ws.Range("E38").Value = "  -8.90%  "
$ws.Range("D39").Value = "'0.0377"
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("E40").Value = "  -1.13%  "
$ws.Range("D41").Value = "'2.55"
$ws.Range("E41").Value = "  -4.01%  "
$ws.Range("D42").Value = "'73.05"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").Value = "'12.63"
$ws.Range("E45").Value = "  -6.92%  "
$ws.Range("D46").Value = "'1.38"
$ws.Range("E46").Value = "  +2.57%  "
$ws.Range("D47").Value = "'5.72"
$ws.Range("E47").Value = "  +2.16%  "
$ws.Range("D48").Value = "'108.12"
$ws.Range("E48").Value = "  +5.83%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'8.84"
$ws.Range("E49").Value = "  +3.38%  "
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "'1.30"
$ws.Range("E50").Value = "  -2.35%  "
$ws.Range("D51").Value = "'73.41"
$ws.Range("E51").Value = "  +4.66%  "
